$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.390.40"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +2.56%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.651.85"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +1.56%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.31%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.48"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.88%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.00"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.79%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.24%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.593"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.74%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.117"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +6.78%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.400"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +3.81%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.81"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +2.34%  "

# Row 12
$ws.Range("E12").Value = "  +1.33%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "29.19"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +4.44%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.121.40"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.16%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "65.144.46"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +2.38%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000173"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +11.63%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.658.83"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.52%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.60"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.36%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.83"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +2.37%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "354.77"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.58%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.27"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +5.63%  "

# Row 22
$ws.Range("E22").Value = "  +0.26%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.98"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.70%  "

# Row 24
$ws.Range("E24").Value = "  +1.09%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.49"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.13%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.67"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.10%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.36"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +4.05%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.165"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.75%  "

# Row 29
$ws.Range("B29").Value = "Bittensor"
$ws.Range("C29").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "545.27"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.71%  "

# Row 30
$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.996"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.40%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0₃0923"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +8.26%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.07"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.66%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.83"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +4.67%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.79"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +9.00%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.46"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +3.86%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.428"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +2.83%  "

# Row 37
$ws.Range("B37").Value = "Stacks"
$ws.Range("C37").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.06"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +5.62%  "

# Row 38
$ws.Range("B38").Value = "Monero"
$ws.Range("C38").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "165.23"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.85%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "20.23"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +3.08%  "

# Row 40
$ws.Range("E40").Value = "  +0.06%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.998"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.15%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "168.59"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.03%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "41.72"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +5.09%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.11"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +3.52%  "

# Row 45
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "23.43"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +5.94%  "

# Row 46
$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0607"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +3.25%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.26"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +11.07%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.648"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.52%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0252"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.04%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0986"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.75%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.52"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.29%  "
